$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 108, shifting existing rows 108-166 down to 109-167.
$ws.Rows(108).Insert()

# Populate the newly inserted row 108 with its data.
$ws.Cells.Item(108, 1).Value = 10
$ws.Cells.Item(108, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(108, 3).Value = "La Araucanía"
$ws.Cells.Item(108, 4).Value = 44806
$ws.Cells.Item(108, 5).Value = 9
$ws.Cells.Item(108, 6).Value = "Fruta"
$ws.Cells.Item(108, 7).Value = 100104
$ws.Cells.Item(108, 8).Value = "Frutos de pepita"
$ws.Cells.Item(108, 9).Value = 100104001
$ws.Cells.Item(108, 10).Value = "Granada"
$ws.Cells.Item(108, 11).Value = "Wonderfull"
$ws.Cells.Item(108, 12).Value = "Primera"
$ws.Cells.Item(108, 13).Value = 160
$ws.Cells.Item(108, 14).Value = 14000
$ws.Cells.Item(108, 15).Value = 15000
$ws.Cells.Item(108, 16).Value = 14500
$ws.Cells.Item(108, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(108, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(108, 19).Value = 1450
$ws.Cells.Item(108, 20).Value = 10
